$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D and E hold numeric/percent-looking text that Excel would
# otherwise auto-convert to numbers; force them to remain plain text by setting
# the NumberFormat to "@" first, and resetting the style back to Normal afterwards
# so no stray style index is left referenced on the cell (matches original which
# has no explicit style on data cells).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '278.81'
Set-TextValue $ws.Range("E2") '6.93%'
Set-TextValue $ws.Range("D3") '27.25'
Set-TextValue $ws.Range("E3") '0.55%'
Set-TextValue $ws.Range("D4") '4.799'
Set-TextValue $ws.Range("E4") '2.09%'
Set-TextValue $ws.Range("D5") '0.06250'
Set-TextValue $ws.Range("E5") '0.41%'
Set-TextValue $ws.Range("E6") '1.81%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D7") '3.263'
Set-TextValue $ws.Range("E7") '2.91%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D8") '0.8774'
Set-TextValue $ws.Range("E8") '2.91%'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range("D9") '0.9446'
Set-TextValue $ws.Range("E9") '3.78%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D10") '0.1442'
Set-TextValue $ws.Range("E10") '2.93%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range("D11") '0.05081'
Set-TextValue $ws.Range("E11") '4.14%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range("D12") '0.07279'
Set-TextValue $ws.Range("E12") '2.74%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D13") '0.03150'
Set-TextValue $ws.Range("E13") '0.76%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D14") '0.09036'
Set-TextValue $ws.Range("E14") '-0.17%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D15") '0.001551'
Set-TextValue $ws.Range("E15") '1.16%'
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range("D16") '0.0006264'
Set-TextValue $ws.Range("E16") '1.81%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D17") '0.005893'
Set-TextValue $ws.Range("E17") '-1.79%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D18") '3.470'
Set-TextValue $ws.Range("E18") '0.52%'
Set-TextValue $ws.Range("E19") '5.06%'
Set-TextValue $ws.Range("E20") '-0.62%'
Set-TextValue $ws.Range("E21") '0.02%'
Set-TextValue $ws.Range("D22") '3.847'
Set-TextValue $ws.Range("E22") '-6.24%'
Set-TextValue $ws.Range("D23") '0.04326'
Set-TextValue $ws.Range("E23") '2.06%'
Set-TextValue $ws.Range("D24") '0.001175'
Set-TextValue $ws.Range("E24") '-3.77%'
Set-TextValue $ws.Range("D25") '0.004275'
Set-TextValue $ws.Range("E25") '4.79%'
Set-TextValue $ws.Range("E26") '-0.20%'
Set-TextValue $ws.Range("E27") '-1.63%'
Set-TextValue $ws.Range("D40") '0.04032'
Set-TextValue $ws.Range("E40") '2.87%'
Set-TextValue $ws.Range("D41") '0.006707'
Set-TextValue $ws.Range("E41") '63.13%'
Set-TextValue $ws.Range("D42") '0.1153'
Set-TextValue $ws.Range("E42") '3.74%'
Set-TextValue $ws.Range("D43") '0.002207'
Set-TextValue $ws.Range("E43") '2.55%'
Set-TextValue $ws.Range("E44") '-9.59%'
Set-TextValue $ws.Range("D45") '0.00005122'
Set-TextValue $ws.Range("E45") '0.15%'
Set-TextValue $ws.Range("E46") '-0.18%'
Set-TextValue $ws.Range("D47") '2.217'
Set-TextValue $ws.Range("E47") '3,086.50%'
Set-TextValue $ws.Range("E49") '-0.18%'
Set-TextValue $ws.Range("E50") '-0.18%'
